# od_aos script corrected; working again
$wb = $excel.ActiveWorkbook

# --- "detailed explanation" sheet (second sheet, index 2) ---
$ws2 = $wb.Worksheets.Item("detailed explanation")

# Insert a new row before row 18, copying formatting from the row above (row 17)
$ws2.Rows("18:18").Insert()

# Correct/complete the previously truncated description (row 17)
$ws2.Range("C17").Value = "folder in locale directory where osmnx output road nodes and edges (wgs84 epsg4326) are located"

# Fill newly inserted row 18 with the pos_source documentation entry
$ws2.Range("B18").Value = "pos_source"
$ws2.Range("C18").Value = "path to source feature for public open space analysis relative to data directory"

# Restore cursor/selection position and make this the active sheet/tab
$ws2.Range("C10").Select()
$ws2.Activate()

# --- "study_regions" sheet ---
$ws4 = $wb.Worksheets.Item("study_regions")
$ws4.Activate()
$ws4.Application.ActiveWindow.ScrollColumn = $ws4.Range("I2").Column
$ws4.Application.ActiveWindow.ScrollRow = $ws4.Range("I2").Row
$ws4.Range("K4").Select()

# Re-activate "detailed explanation" so it remains the active/selected tab
$ws2.Activate()
